$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add I1 = "I0" and J1 = "IF", matching the formatting of the existing header cells ---
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows 2..31: add values for columns I (I0) and J (IF) ---
$data = @{
    2  = @(2, 2)
    3  = @(3, 3)
    4  = @(9, 9)
    5  = @(1, 2)
    6  = @(8, 8)
    7  = @(10, 10)
    8  = @(1, 3)
    9  = @(9, 9)
    10 = @(8, 9)
    11 = @(9, 9)
    12 = @(7, 7)
    13 = @(8, 8)
    14 = @(9, 9)
    15 = @(8, 8)
    16 = @(5, 5)
    17 = @(7, 7)
    18 = @(8, 8)
    19 = @(7, 7)
    20 = @(8, 8)
    21 = @(9, 9)
    22 = @(6, 7)
    23 = @(7, 7)
    24 = @(6, 6)
    25 = @(7, 8)
    26 = @(7, 7)
    27 = @(8, 8)
    28 = @(5, 5)
    29 = @(9, 9)
    30 = @(5, 6)
    31 = @(3, 3)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
